$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 5.5
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 2.35
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 3.75
$ws.Range("G4").Value = 3.5
$ws.Range("H4").Value = 2.62
$ws.Range("J4").Value = 4.2
$ws.Range("K4").Value = 1.8
$ws.Range("L4").Value = 3.2
$ws.Range("M4").Value = 1.15
$ws.Range("N4").Value = 4.7
$ws.Range("Q4").Value = 2.8
$ws.Range("W4").Value = 1.62
$ws.Range("X4").Value = 2.15
$ws.Range("Y4").Value = 2.2
$ws.Range("AE4").Value = 40
$ws.Range("AG4").Value = 4.7
$ws.Range("AH4").Value = 5.4
$ws.Range("AN4").Value = 10.25
$ws.Range("AO4").Value = 26
$ws.Range("G5").Value = 4.25
$ws.Range("H5").Value = 2.8
$ws.Range("I5").Value = 2.05
$ws.Range("J5").Value = 4.7
$ws.Range("K5").Value = 1.87
$ws.Range("L5").Value = 2.75
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 5.1
$ws.Range("O5").Value = 1.55
$ws.Range("P5").Value = 2.32
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.44
$ws.Range("U5").Value = 4.65
$ws.Range("W5").Value = 1.57
$ws.Range("X5").Value = 2.25
$ws.Range("Y5").Value = 2.15
$ws.Range("Z5").Value = 1.62
$ws.Range("AA5").Value = 8.75
$ws.Range("AB5").Value = 22
$ws.Range("AC5").Value = 14.5
$ws.Range("AD5").Value = 75
$ws.Range("AE5").Value = 50
$ws.Range("AF5").Value = 65
$ws.Range("AG5").Value = 5.1
$ws.Range("AH5").Value = 5.7
$ws.Range("AI5").Value = 18
$ws.Range("AL5").Value = 5.2
$ws.Range("AM5").Value = 8.25
$ws.Range("AN5").Value = 9.25
$ws.Range("AO5").Value = 19
$ws.Range("AP5").Value = 21
$ws.Range("G6").Value = 1.65
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 5.5
$ws.Range("J6").Value = 2.2
$ws.Range("K6").Value = 2.12
$ws.Range("L6").Value = 5.5
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 6.4
$ws.Range("O6").Value = 1.38
$ws.Range("P6").Value = 2.82
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.65
$ws.Range("U6").Value = 3.6
$ws.Range("V6").Value = 1.25
$ws.Range("W6").Value = 1.42
$ws.Range("X6").Value = 2.67
$ws.Range("Y6").Value = 2.02
$ws.Range("AA6").Value = 5.7
$ws.Range("AB6").Value = 7
$ws.Range("AD6").Value = 12.5
$ws.Range("AE6").Value = 14.5
$ws.Range("AG6").Value = 6.4
$ws.Range("AH6").Value = 6.6
$ws.Range("AI6").Value = 18
$ws.Range("AL6").Value = 12
$ws.Range("AM6").Value = 30
$ws.Range("AO6").Value = 110
$ws.Range("AQ6").Value = 65
$ws.Range("Q7").Value = 2.05
$ws.Range("R7").Value = 1.68
$ws.Range("Q8").Value = 1.98
$ws.Range("R8").Value = 1.88
$ws.Range("AD8").Value = 7.5
$ws.Range("AG8").Value = 9.5
$ws.Range("AH8").Value = 10
$ws.Range("AI8").Value = 29
$ws.Range("AJ8").Value = 101
$ws.Range("AL8").Value = 19
$ws.Range("H11").Value = 3.25
$ws.Range("Q11").Value = 2.05
$ws.Range("R11").Value = 1.68
$ws.Range("Y11").Value = 1.8
$ws.Range("Z11").Value = 1.8
$ws.Range("Y12").Value = 1.69
$ws.Range("Q13").Value = 1.72
$ws.Range("W13").Value = 1.36
$ws.Range("Y13").Value = 1.63
$ws.Range("Q14").Value = 1.47
$ws.Range("W14").Value = 1.25
$ws.Range("Y14").Value = 1.58
$ws.Range("Q15").Value = 1.72
$ws.Range("W15").Value = 1.33
$ws.Range("Y15").Value = 1.54
$ws.Range("W16").Value = 1.36
$ws.Range("Y16").Value = 1.63
$ws.Range("I17").Value = 1.69
$ws.Range("W17").Value = 1.3
$ws.Range("Y17").Value = 1.58
$ws.Range("G18").Value = 1.47
$ws.Range("H18").Value = 3.85
$ws.Range("I18").Value = 6.6
$ws.Range("J18").Value = 1.98
$ws.Range("L18").Value = 6.2
$ws.Range("O18").Value = 1.31
$ws.Range("P18").Value = 2.87
$ws.Range("Q18").Value = 1.93
$ws.Range("R18").Value = 1.7
$ws.Range("U18").Value = 3.15
$ws.Range("V18").Value = 1.26
$ws.Range("AB18").Value = 6.1
$ws.Range("AD18").Value = 9.5
$ws.Range("AE18").Value = 13.5
$ws.Range("AG18").Value = 9
$ws.Range("AH18").Value = 7.7
$ws.Range("AI18").Value = 22
$ws.Range("AL18").Value = 15
$ws.Range("AM18").Value = 40
$ws.Range("AN18").Value = 21
$ws.Range("AO18").Value = 150
$ws.Range("AP18").Value = 80
$ws.Range("AQ18").Value = 80
$ws.Range("G19").Value = 1.5
$ws.Range("H19").Value = 3.9
$ws.Range("I19").Value = 6.1
$ws.Range("L19").Value = 5.8
$ws.Range("R19").Value = 1.9
$ws.Range("Y19").Value = 1.82
$ws.Range("Z19").Value = 1.78
$ws.Range("AA19").Value = 6.7
$ws.Range("AB19").Value = 7
$ws.Range("AD19").Value = 10.5
$ws.Range("AF19").Value = 26
$ws.Range("AH19").Value = 7.7
$ws.Range("AL19").Value = 16.5
$ws.Range("AM19").Value = 40
$ws.Range("AN19").Value = 19
$ws.Range("AO19").Value = 150
$ws.Range("G26").Value = 1.5
$ws.Range("H26").Value = 3.85
$ws.Range("I26").Value = 6.4
$ws.Range("J26").Value = 2.07
$ws.Range("K26").Value = 2.15
$ws.Range("L26").Value = 6.1
$ws.Range("M26").Value = 1.06
$ws.Range("N26").Value = 7.3
$ws.Range("O26").Value = 1.3
$ws.Range("P26").Value = 3.25
$ws.Range("Q26").Value = 1.88
$ws.Range("R26").Value = 1.82
$ws.Range("U26").Value = 3.1
$ws.Range("V26").Value = 1.32
$ws.Range("W26").Value = 1.42
$ws.Range("X26").Value = 2.65
$ws.Range("Y26").Value = 2
$ws.Range("Z26").Value = 1.72
$ws.Range("AA26").Value = 6.1
$ws.Range("AB26").Value = 6.6
$ws.Range("AC26").Value = 8.25
$ws.Range("AD26").Value = 10.25
$ws.Range("AE26").Value = 12.5
$ws.Range("AF26").Value = 29
$ws.Range("AG26").Value = 7.3
$ws.Range("AH26").Value = 7.6
$ws.Range("AI26").Value = 18.5
$ws.Range("AJ26").Value = 100
$ws.Range("AK26").Value = 800
$ws.Range("AL26").Value = 15.5
$ws.Range("AM26").Value = 40
$ws.Range("AN26").Value = 20
$ws.Range("AO26").Value = 150
$ws.Range("AP26").Value = 75
$ws.Range("AQ26").Value = 70
